# Add "hybrid bold + color" highlighting to the quantitative impact metrics
# (percentages, dollar amounts, large counts) called out in specific resume
# bullet points. For each target bullet we:
#   1. Find.Execute a long, unique anchor string (the whole bullet sentence,
#      or a unique sub-phrase of it) so we never touch a different paragraph
#      that happens to contain the same bare number (e.g. "23%" shows up in
#      the summary, the key-projects section, *and* the bullet we want).
#   2. Within that matched Range, locate the metric substring(s) by text
#      offset and build a sub-Range for just that substring.
#   3. Set Font.Bold / Font.Color on the sub-Range, which causes Word to
#      split the run exactly the way native Word COM editing would.

$d = $word.ActiveDocument

# Word's Font.Color is an OLE_COLOR long encoded 0x00BBGGRR. The target
# color is #2C3E50 (R=2C G=3E B=50).
$highlightColor = 0x2C + (0x3E * 256) + (0x50 * 65536)   # 5258796

function Set-MetricBold([object]$range, [string]$metric) {
    # range: the outer Range already positioned over an anchor phrase.
    # metric: the exact substring inside that anchor's text to bold+color.
    $fullText = $range.Text
    $idx = $fullText.IndexOf($metric)
    if ($idx -lt 0) {
        throw "Metric '$metric' not found inside anchor text '$fullText'"
    }
    $subStart = $range.Start + $idx
    $subEnd = $subStart + $metric.Length
    $sub = $d.Range($subStart, $subEnd)
    $sub.Font.Bold = 1
    $sub.Font.Color = $highlightColor
}

function Find-Anchor([string]$anchor) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor not found: $anchor"
    }
    return $rng
}

# 1. Siege Analytics bullet - demographic classification accuracy
$anchor = Find-Anchor "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%"
Set-MetricBold $anchor "23%"
Set-MetricBold $anchor "64%"

# 2. Siege Analytics bullet - survey margin of error / turnout prediction
$anchor = Find-Anchor "Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes"
Set-MetricBold $anchor "±4.2%"
Set-MetricBold $anchor "±2.1%"
Set-MetricBold $anchor "71%"
Set-MetricBold $anchor "87%"

# 3. Siege Analytics bullet - boundary estimation / mapping costs
$anchor = Find-Anchor "Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis"
Set-MetricBold $anchor "73.5%"
Set-MetricBold $anchor "`$4.7M"

# 4. Siege Analytics bullet - FEC analysis / political spending sub-economy
$anchor = Find-Anchor "Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion"
Set-MetricBold $anchor "`$2"

# 5. Helm/Murmuration bullet - ETL modernization processing time
$anchor = Find-Anchor "Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%"
Set-MetricBold $anchor "57%"

# 6. Key Achievements - platform impact / analysts served
$anchor = Find-Anchor "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
Set-MetricBold $anchor "12,847"

# 7. Key Achievements - revenue generation
$anchor = Find-Anchor "Revenue generation: Delivered `$4.9M additional revenue through optimization"
Set-MetricBold $anchor "`$4.9M"

# 8. Key Achievements - conversion rate improvement
$anchor = Find-Anchor "23% conversion rate improvement"
Set-MetricBold $anchor "23%"
